# The commit swaps the contents of ppt/theme/theme1.xml ("Integral" / "Red
# Violet" colour scheme, used by the single slide master -> all slides) and
# ppt/theme/theme2.xml ("Office Theme" / "Office" colour scheme, used only by
# the notes master). Net visible effect on the deck: the slide master's
# theme colours change from the "Integral" palette to the stock "Office"
# palette. Font scheme and format scheme (fills/lines/effects) are identical
# between the two themes, so only the 12 colour-scheme slots need updating.
#
# COM-interop RGB values are packed 0x00BBGGRR (reverse byte order of the
# "RRGGBB" hex used in the OOXML <a:srgbClr val="RRGGBB"/>).
#
# Target ("Office") theme colours, in ThemeColorScheme.Item(1..12) order
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink):
#   1  dk1       000000 -> 0
#   2  lt1       FFFFFF -> 16777215
#   3  dk2       44546A -> 6968388
#   4  lt2       E7E6E6 -> 15132391
#   5  accent1   5B9BD5 -> 13998939
#   6  accent2   ED7D31 -> 3243501
#   7  accent3   A5A5A5 -> 10855845
#   8  accent4   FFC000 -> 49407
#   9  accent5   4472C4 -> 12874308
#  10  accent6   70AD47 -> 4697456
#  11  hlink     0563C1 -> 12673797
#  12  folHlink  954F72 -> 7491477

$p = $ppt.ActivePresentation

$officeThemeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

$slide = $p.Slides.Item(1)
$tcs = $slide.ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $officeThemeColors[$i - 1]
}
